# Update the "Förändrad" (Changed) date column (C) for rows 2-13
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05),
# matching the data update recorded in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
